# Rename the three logo InlineShapes that live in the document's
# header/footer stories. The diff only touches the DrawingML "name"
# label (wp:docPr / pic:cNvPr) that Word stamps on a pasted picture -
# the picture content, size, and position are unchanged.
#
#   footer (default)    - Pearson logo: image1.png -> image2.png
#   footer (first page) - Pearson logo: image1.png -> image2.png
#   header (first page) - BTEC logo:    image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$section = $d.Sections(1)

# Default footer - Pearson Edexcel logo (docPr id="1").
$footerDefault = $section.Footers(1)
$footerDefault.Range.InlineShapes(1).Name = "image2.png"

# First-page footer - Pearson Edexcel logo (docPr id="2").
$footerFirst = $section.Footers(2)
$footerFirst.Range.InlineShapes(1).Name = "image2.png"

# First-page header - BTEC logo (docPr id="3").
$headerFirst = $section.Headers(2)
$headerFirst.Range.InlineShapes(1).Name = "image1.jpg"
